$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet/tab
$ws.Name = "Ravichandran Ashwin"

# Force the numeric-looking stat columns (runs, balls, fours, sixes, sr) to be
# stored as text so values like "112.50" / "100.00" keep their exact string
# representation instead of being coerced into numbers.
$ws.Range("E2:I6").NumberFormat = "@"

# Header row (row 1) - insert "matchNo" as new column A, shifting everything right
$ws.Range("A1").Value = "matchNo"
$ws.Range("B1").Value = "teamName"
$ws.Range("C1").Value = "batterName"
$ws.Range("D1").Value = "states"
$ws.Range("E1").Value = "runs"
$ws.Range("F1").Value = "balls"
$ws.Range("G1").Value = "fours"
$ws.Range("H1").Value = "sixes"
$ws.Range("I1").Value = "sr"
$ws.Range("J1").Value = "opponentTeamName"
$ws.Range("K1").Value = "venue"
$ws.Range("L1").Value = "date"
$ws.Range("M1").Value = "result"

# Row 2
$ws.Range("A2").Value = "41st"
$ws.Range("B2").Value = "Delhi Capitals"
$ws.Range("C2").Value = "Ravichandran Ashwin"
$ws.Range("D2").Value = "c Rana b Southee"
$ws.Range("E2").Value = "9"
$ws.Range("F2").Value = "8"
$ws.Range("G2").Value = "1"
$ws.Range("H2").Value = "0"
$ws.Range("I2").Value = "112.50"
$ws.Range("J2").Value = "Kolkata Knight Riders"
$ws.Range("K2").Value = "Sharjah"
$ws.Range("L2").Value = "September 28"
$ws.Range("M2").Value = "KKR won by 3 wickets (with 10 balls remaining)"

# Row 3
$ws.Range("A3").Value = "50th"
$ws.Range("B3").Value = "Delhi Capitals"
$ws.Range("C3").Value = "Ravichandran Ashwin"
$ws.Range("D3").Value = "b Thakur"
$ws.Range("E3").Value = "2"
$ws.Range("F3").Value = "3"
$ws.Range("G3").Value = "0"
$ws.Range("H3").Value = "0"
$ws.Range("I3").Value = "66.66"
$ws.Range("J3").Value = "Chennai Super Kings"
$ws.Range("K3").Value = "Dubai (DSC)"
$ws.Range("L3").Value = "October 04"
$ws.Range("M3").Value = "Capitals won by 3 wickets (with 2 balls remaining)"

# Row 4 (states/D4 is blank in the source data)
$ws.Range("A4").Value = "36th"
$ws.Range("B4").Value = "Delhi Capitals"
$ws.Range("C4").Value = "Ravichandran Ashwin"
$ws.Range("E4").Value = "6"
$ws.Range("F4").Value = "6"
$ws.Range("G4").Value = "0"
$ws.Range("H4").Value = "0"
$ws.Range("I4").Value = "100.00"
$ws.Range("J4").Value = "Rajasthan Royals"
$ws.Range("K4").Value = "Abu Dhabi"
$ws.Range("L4").Value = "September 25"
$ws.Range("M4").Value = "Capitals won by 33 runs"

# Row 5 (states/D5 is blank in the source data)
$ws.Range("A5").Value = "46th"
$ws.Range("B5").Value = "Delhi Capitals"
$ws.Range("C5").Value = "Ravichandran Ashwin"
$ws.Range("E5").Value = "20"
$ws.Range("F5").Value = "21"
$ws.Range("G5").Value = "0"
$ws.Range("H5").Value = "1"
$ws.Range("I5").Value = "95.23"
$ws.Range("J5").Value = "Mumbai Indians"
$ws.Range("K5").Value = "Sharjah"
$ws.Range("L5").Value = "October 02"
$ws.Range("M5").Value = "Capitals won by 4 wickets (with 5 balls remaining)"

# Row 6 (original row 2 data, shifted right by one column)
$ws.Range("A6").Value = "7th"
$ws.Range("B6").Value = "Delhi Capitals"
$ws.Range("C6").Value = "Ravichandran Ashwin"
$ws.Range("D6").Value = "run out (Miller/†Samson)"
$ws.Range("E6").Value = "7"
$ws.Range("F6").Value = "4"
$ws.Range("G6").Value = "1"
$ws.Range("H6").Value = "0"
$ws.Range("I6").Value = "175.00"
$ws.Range("J6").Value = "Rajasthan Royals"
$ws.Range("K6").Value = "Wankhede"
$ws.Range("L6").Value = "April 15"
$ws.Range("M6").Value = "Royals won by 3 wickets (with 2 balls remaining)"
